$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 227-229 continue the daily series in column A, which uses the
# same date/time style (s="2") as every other row above it. Copy that
# formatting down before writing the new values.
$xlPasteFormats = -4122
$ws.Range("A226").Copy()
$ws.Range("A227:A229").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 0
$ws.Range("C227").Value = 8
$ws.Range("D227").Value = 250.0781494216943

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 1
$ws.Range("C228").Value = 7
$ws.Range("D228").Value = 218.8183807439825

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 0
$ws.Range("C229").Value = 5
$ws.Range("D229").Value = 156.2988433885589
